# Realestate Update resale numbers 2024-01-10 12:51
# Append a new data row (row 42) to the CityResaleNum sheet with the
# latest resale-number snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

# --- Text columns -----------------------------------------------------
# These values look numeric/date-like ("2024-01-10", "01") so Excel would
# normally auto-convert them on assignment (date serial / drop leading
# zero). Force text storage via NumberFormat "@" then clear the
# formatting again afterwards so the cell ends up with the default style
# (matching the rest of the sheet) while keeping the text value.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-10"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "12:51:32"

$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "01"
$ws.Cells.Item($row, 4).ClearFormats()

# --- Numeric columns ----------------------------------------------------
$ws.Cells.Item($row, 5).Value = 139521
$ws.Cells.Item($row, 6).Value = 142698
$ws.Cells.Item($row, 7).Value = 171925
$ws.Cells.Item($row, 8).Value = 147843
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118893
$ws.Cells.Item($row, 11).Value = 224703
$ws.Cells.Item($row, 12).Value = 250839
$ws.Cells.Item($row, 13).Value = 185057
$ws.Cells.Item($row, 14).Value = 110484
$ws.Cells.Item($row, 15).Value = 40721
$ws.Cells.Item($row, 16).Value = 30876
$ws.Cells.Item($row, 17).Value = 72659
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 41796
$ws.Cells.Item($row, 20).Value = -1
